# Fit of clock controlled IAV infection model.
# Applies: Sheet2 viral-titer row restructure + new Sheet3 (per-ZT breakdown)
# plus the associated selection / active-sheet bookkeeping.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet2: row 8 (time axis for the "Viral titers" block) gains 4 more time
# points (24/48/96/144) in B:E, pushing the existing 8*24 / 10*24 formulas
# from B8:C8 out to F8:G8 (keeping their bold/bordered style).
# ---------------------------------------------------------------------------
$ws2.Range("B8:C8").Copy()
$ws2.Range("F8").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Range("F8").Formula = "=8*24"
$ws2.Range("G8").Formula = "=10*24"

$ws2.Range("D9").Copy()                # plain numeric style (s=31) donor
$ws2.Range("B8").PasteSpecial(-4122)
$ws2.Range("C8").PasteSpecial(-4122)
$ws2.Range("D8").PasteSpecial(-4122)
$ws2.Range("E8").PasteSpecial(-4122)

$ws2.Range("B8").Value = 24
$ws2.Range("C8").Value = 48
$ws2.Range("D8").Value = 96
$ws2.Range("E8").Value = 144

# Row 9 ("ZT23"): new values at the 4 extra time points, old values moved to F:G
$ws2.Range("B9").Value = 0.05
$ws2.Range("C9").Value = 0.08
$ws2.Range("D9").Value = 0.1
$ws2.Range("E9").Value = 0.5
$ws2.Range("F9").Value = 3.0320002526666663
$ws2.Range("G9").Value = 3.6068164140000007

# Row 10 ("ZT11"): same treatment
$ws2.Range("B10").Value = 0.05
$ws2.Range("C10").Value = 0.1
$ws2.Range("D10").Value = 0.5
$ws2.Range("E10").Value = 1
$ws2.Range("F10").Value = 14.149333333333335
$ws2.Range("G10").Value = 19.133183586000001

# ---------------------------------------------------------------------------
# View bookkeeping: Sheet1's selection moves to B90 (scroll pos untouched,
# cosmetic), Sheet2's selection widens from B22 to B22:G22, and the tab
# selection ultimately ends on the new Sheet3 (handled below). Sheet1/2 are
# selected first so the later Sheet3 selection is what "wins" tabSelected.
# ---------------------------------------------------------------------------
$ws1.Range("B90").Select()
$ws2.Range("B22:G22").Select()

# ---------------------------------------------------------------------------
# New Sheet3: per-ZT (ZT23/ZT11) breakdown of M, NK, V (viral titer) and the
# two T-cell series (T / T_E) used by the clock-controlled IAV infection fit.
# Added after Sheet2 so it lands last / becomes the active tab.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"

# Block 1: M (macrophages?) ZT23 vs ZT11
$ws3.Range("A1").Value = "M"
$ws3.Range("B1").Value = 24
$ws3.Range("C1").Value = 48
$ws3.Range("D1").Value = 96
$ws3.Range("E1").Value = 144

$ws3.Range("A2").Value = "ZT23"
$ws3.Range("B2").Value = 13.02786
$ws3.Range("C2").Value = 17.766539999999999
$ws3.Range("D2").Value = 23.985200000000003
$ws3.Range("E2").Value = 40.158839999999998

$ws3.Range("A3").Value = "ZT11"
$ws3.Range("B3").Value = 19.958400000000001
$ws3.Range("C3").Value = 28.712300000000003
$ws3.Range("D3").Value = 36.297600000000003
$ws3.Range("E3").Value = 33.307200000000002

# Block 2: NK cells ZT23 vs ZT11
$ws3.Range("A5").Value = "NK"
$ws3.Range("B5").Value = 24
$ws3.Range("C5").Value = 48
$ws3.Range("D5").Value = 96
$ws3.Range("E5").Value = 144

$ws3.Range("A6").Value = "ZT23"
$ws3.Range("B6").Value = 5.2554000000000007
$ws3.Range("C6").Value = 6.7003199999999996
$ws3.Range("D6").Value = 7.9910000000000005
$ws3.Range("E6").Value = 9.3234000000000012

$ws3.Range("A7").Value = "ZT11"
$ws3.Range("B7").Value = 6.2207999999999997
$ws3.Range("C7").Value = 4.8201499999999999
$ws3.Range("D7").Value = 7.1440000000000001
$ws3.Range("E7").Value = 7.7759999999999989

# Block 3: V (viral titers) ZT23 vs ZT11, 8 time points, mirrors Sheet2 row 5/6
$ws3.Range("A9").Value = "V"
$ws3.Range("B9").Value = 6
$ws3.Range("C9").Value = 12
$ws3.Range("D9").Value = 24
$ws3.Range("E9").Value = 48
$ws3.Range("F9").Value = 96
$ws3.Range("G9").Value = 144
$ws3.Range("H9").Formula = "=8*24"
$ws3.Range("I9").Formula = "=10*24"

$ws3.Range("A10").Value = "ZT23"
$ws3.Range("B10").Value = 0
$ws3.Range("C10").Value = 1.65
$ws3.Range("D10").Value = 4.6000000000000005
$ws3.Range("E10").Value = 4.7250000000000005
$ws3.Range("F10").Value = 5.55
$ws3.Range("G10").Value = 5.1000000000000005
$ws3.Range("H10").Value = 3.4249999999999998
$ws3.Range("I10").Value = 0

$ws3.Range("A11").Value = "ZT11"
$ws3.Range("B11").Value = 0
$ws3.Range("C11").Value = 2.0249999999999999
$ws3.Range("D11").Value = 5.1749999999999998
$ws3.Range("E11").Value = 4.7750000000000004
$ws3.Range("F11").Value = 5.6749999999999998
$ws3.Range("G11").Value = 4.6500000000000004
$ws3.Range("H11").Value = 4.5
$ws3.Range("I11").Value = 0.35

# Give B9:C9 / H9:I9 the bold/bordered header style (s=39), matching Sheet2's
# analogous "time axis" row 5 (B5:C5 / H5:I5 are s=39 there too).
$ws2.Range("B5").Copy()
$ws3.Range("B9").PasteSpecial(-4122)
$ws3.Range("C9").PasteSpecial(-4122)
$ws3.Range("H9").PasteSpecial(-4122)
$ws3.Range("I9").PasteSpecial(-4122)

# Block 4: T ZT23 vs ZT11, 6 time points + two formula columns
$ws3.Range("A13").Value = "T"
$ws3.Range("B13").Value = 24
$ws3.Range("C13").Value = 48
$ws3.Range("D13").Value = 96
$ws3.Range("E13").Value = 144
$ws3.Range("F13").Formula = "=8*24"
$ws3.Range("G13").Formula = "=10*24"

$ws3.Range("A14").Value = "ZT23"
$ws3.Range("B14").Value = 0.05
$ws3.Range("C14").Value = 0.08
$ws3.Range("D14").Value = 0.1
$ws3.Range("E14").Value = 0.5
$ws3.Range("F14").Value = 3.0320002526666663
$ws3.Range("G14").Value = 3.6068164140000007

$ws3.Range("A15").Value = "ZT11"
$ws3.Range("B15").Value = 0.05
$ws3.Range("C15").Value = 0.08
$ws3.Range("D15").Value = 0.1
$ws3.Range("E15").Value = 0.5
$ws3.Range("F15").Value = 2.3328000000000002
$ws3.Range("G15").Value = 3.4830000000000001

# Block 5: T_E ZT23 vs ZT11, 6 time points + two formula columns
$ws3.Range("A17").Value = "T_E"
$ws3.Range("B17").Value = 24
$ws3.Range("C17").Value = 48
$ws3.Range("D17").Value = 96
$ws3.Range("E17").Value = 144
$ws3.Range("F17").Formula = "=8*24"
$ws3.Range("G17").Formula = "=10*24"

$ws3.Range("A18").Value = "ZT23"
$ws3.Range("B18").Value = 0.05
$ws3.Range("C18").Value = 0.1
$ws3.Range("D18").Value = 0.5
$ws3.Range("E18").Value = 1
$ws3.Range("F18").Value = 14.149333333333335
$ws3.Range("G18").Value = 19.133183586000001

$ws3.Range("A19").Value = "ZT11"
$ws3.Range("B19").Value = 0.05
$ws3.Range("C19").Value = 0.1
$ws3.Range("D19").Value = 0.5
$ws3.Range("E19").Value = 1
$ws3.Range("F19").Value = 10.160639999999999
$ws3.Range("G19").Value = 15.1632

# Style rows 13 and 17's formula cells (F/G) bold+bordered like row 9's H/I.
$ws2.Range("H5").Copy()
$ws3.Range("F13").PasteSpecial(-4122)
$ws3.Range("G13").PasteSpecial(-4122)
$ws3.Range("F17").PasteSpecial(-4122)
$ws3.Range("G17").PasteSpecial(-4122)

# Final selection/view state: Sheet3 tab is active, cursor at C8.
$ws3.Range("C8").Select()

Write-Output "done"
